$wb = $excel.ActiveWorkbook

# Map of sheet name -> new values for B2:E2 (Time taken, distance offset, distance traveled, area between paths)
$updates = @{
    "Trial No.1"  = @(64.801, 0.1585722829283145, 33.73754501342773, 3.838452816009521)
    "Trial No.2"  = @(64.203, 0.1426898343830232, 33.20681381225586, 3.620044708251953)
    "Trial No.3"  = @(66.19499999999999, 0.104727913956677, 33.45143890380859, 2.84593391418457)
    "Trial No.4"  = @(68.803, 0.1438927158829818, 33.9770393371582, 3.392913341522217)
    "Trial No.5"  = @(63.024, 0.0952277880090163, 33.0434684753418, 2.709926128387451)
    "Trial No.6"  = @(66.199, 0.1817085484422198, 33.88814926147461, 3.349871873855591)
    "Trial No.7"  = @(66.607, 0.09050246661430196, 33.84307098388672, 3.04830002784729)
    "Trial No.8"  = @(65.80199999999999, 0.1211052528498621, 33.83699798583984, 3.128567934036255)
    "Trial No.9"  = @(64.006, 0.1328935128820808, 33.19462585449219, 2.888601779937744)
    "Trial No.10" = @(63.607, 0.150469518016676, 33.10587310791016, 3.36250901222229)
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $updates[$sheetName]
    $ws.Range("B2").Value = $vals[0]
    $ws.Range("C2").Value = $vals[1]
    $ws.Range("D2").Value = $vals[2]
    $ws.Range("E2").Value = $vals[3]
}
